$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# Note: ColumnWidth round-trips through the stored XML `width` attribute with
# a constant +0.833333 (5/6) offset in this engine, so we pre-compensate by
# subtracting 5/6 from each desired target width.
$ws.Columns.Item(2).ColumnWidth = 55 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 80 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 57 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 17 - (5/6)
$ws.Columns.Item(7).ColumnWidth = 16 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 63 - (5/6)

# --- New row data (rows 2-10) ---
$rows = @(
    @("1330558", "https://aiesec.org/opportunity/global-talent/1330558", "Technical Account Manager", "Amman, Jordan", "No", "15 applicants", "9 - 12 Weeks", "Barq"),
    @("1330306", "https://aiesec.org/opportunity/global-talent/1330306", "Neuro-Marketing & Communications Intern", "Amman, Jordan", "No", "39 applicants", "9 - 12 Weeks", "Amoux Group"),
    @("1330301", "https://aiesec.org/opportunity/global-talent/1330301", "Business Development Officer", "Amman, Jordan", "No", "35 applicants", "9 - 12 Weeks", "International TEFL Training Institute"),
    @("1330065", "https://aiesec.org/opportunity/global-talent/1330065", "[EXP] Purchase to Pay Process in a Global Context (Spanish Speaker Preferred)", "Maastricht, Netherlands", "Yes", "98 applicants", "6 - 18 Months", "DHL Group"),
    @("1328490", "https://aiesec.org/opportunity/global-talent/1328490", "Sales Intern", "Ümraniye, Elmalıkent, 34764 Ümraniye/İstanbul, Türkiye", "No", "76 applicants", "9 - 12 Weeks", "ENTES ELEKTRONİK CİHAZLAR İMALAT VE TİCARET ANONİM ŞİRKETİ"),
    @("1327381", "https://aiesec.org/opportunity/global-talent/1327381", "Product Management Intern", "Ümraniye, Elmalıkent, 34764 Ümraniye/İstanbul, Türkiye", "No", "113 applicants", "9 - 12 Weeks", "ENTES ELEKTRONİK CİHAZLAR İMALAT VE TİCARET ANONİM ŞİRKETİ"),
    @("1327380", "https://aiesec.org/opportunity/global-talent/1327380", "Comunication Intern", "Ümraniye, Elmalıkent, 34764 Ümraniye/İstanbul, Türkiye", "No", "105 applicants", "9 - 12 Weeks", "ENTES ELEKTRONİK CİHAZLAR İMALAT VE TİCARET ANONİM ŞİRKETİ"),
    @("1321497", "https://aiesec.org/opportunity/global-talent/1321497", "Sales  Specialist", "Kartepe, Kocaeli, Türkiye", "No", "71 applicants", "6 - 18 Months", "Dessa Teknoloji Sanayi Ticaret Limited Şirketi"),
    @("1320725", "https://aiesec.org/opportunity/global-talent/1320725", "International Educational Consultant", "İstanbul, Türkiye", "No", "52 applicants", "6 - 18 Months", "JOHN AND JOHN EĞİTİM TEKNOLOJİ VE İNTERNET YATIRIMLARI LİMİT")
)

$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

# --- Highlight E5 with yellow fill (Yes / premium) ---
$e5 = $ws.Range("E5")
$e5.Interior.Color = 65535
